$p = $ppt.ActivePresentation

# --- Update the auto-date placeholder text ("datetimeFigureOut" fields)
# on the slide master and every slide layout from 1/25/21 -> 7/27/21.
$newDate = "7/27/21"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $m.CustomLayouts.Count; $L++) {
    $lay = $m.CustomLayouts.Item($L)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $shp = $lay.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Fix the typo on the "Advantages" slide title (slide 4).
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "ADVANTAGES"
